$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "JRmMF167"
$ws.Range("B2").Value = 23092933
$ws.Range("C2").Value = "pcntnug55"
$ws.Range("D2").Value = "t&RM9!x6"
$ws.Range("F2").Value = "zdJMQuUb"
$ws.Range("G2").Value = "fOTL"
